$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 20.030895
$ws.Range("H2").Value = 60.092685
$ws.Range("I2").Value = 0.1333691355055119
$ws.Range("J2").Value = 0.1333691355055119
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.580817
$ws.Range("N2").Value = 1.742451
$ws.Range("O2").Value = 0.1705489461834183
$ws.Range("P2").Value = 0.1705489461834183
$ws.Range("Q2").Value = 11.634284341215
$ws.Range("R2").Value = 104.708559070935
$ws.Range("S2").Value = 0.02274596551385858
$ws.Range("T2").Value = 0.02274596551385858

# Row 3
$ws.Range("G3").Value = 20.030895
$ws.Range("H3").Value = 60.092685
$ws.Range("I3").Value = 0.1333691355055119
$ws.Range("J3").Value = 0.1333691355055119
$ws.Range("O3").Value = 0.6679715536912479
$ws.Range("P3").Value = 0.6679715536912479
$ws.Range("Q3").Value = 45.566807426235
$ws.Range("R3").Value = 410.101266836115
$ws.Range("S3").Value = 0.08908678865807539
$ws.Range("T3").Value = 0.08908678865807539

# Row 4
$ws.Range("G4").Value = 20.030895
$ws.Range("H4").Value = 60.092685
$ws.Range("I4").Value = 0.1333691355055119
$ws.Range("J4").Value = 0.1333691355055119
$ws.Range("M4").Value = 0.5499303333333333
$ws.Range("N4").Value = 1.649791
$ws.Range("O4").Value = 0.1614795001253337
$ws.Range("P4").Value = 0.1614795001253338
$ws.Range("Q4").Value = 11.015596764315
$ws.Range("R4").Value = 99.140370878835
$ws.Range("S4").Value = 0.02153638133357797
$ws.Range("T4").Value = 0.02153638133357797

# Row 5
$ws.Range("I5").Value = 0.6531407302146811
$ws.Range("J5").Value = 0.653140730214681
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.580817
$ws.Range("N5").Value = 1.742451
$ws.Range("O5").Value = 0.1705489461834183
$ws.Range("P5").Value = 0.1705489461834183
$ws.Range("Q5").Value = 56.97588832187
$ws.Range("R5").Value = 512.7829948968299
$ws.Range("S5").Value = 0.1113924632475822
$ws.Range("T5").Value = 0.1113924632475822

# Row 6
$ws.Range("I6").Value = 0.6531407302146811
$ws.Range("J6").Value = 0.653140730214681
$ws.Range("O6").Value = 0.6679715536912479
$ws.Range("P6").Value = 0.6679715536912479
$ws.Range("S6").Value = 0.4362794283405367
$ws.Range("T6").Value = 0.4362794283405367

# Row 7
$ws.Range("I7").Value = 0.6531407302146811
$ws.Range("J7").Value = 0.653140730214681
$ws.Range("M7").Value = 0.5499303333333333
$ws.Range("N7").Value = 1.649791
$ws.Range("O7").Value = 0.1614795001253337
$ws.Range("P7").Value = 0.1614795001253338
$ws.Range("Q7").Value = 53.94602647100333
$ws.Range("R7").Value = 485.51423823903
$ws.Range("S7").Value = 0.1054688386265622
$ws.Range("T7").Value = 0.1054688386265622

# Row 8
$ws.Range("G8").Value = 3.916733333333333
$ws.Range("H8").Value = 11.7502
$ws.Range("I8").Value = 0.0260782825067122
$ws.Range("J8").Value = 0.02607828250671219
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.580817
$ws.Range("N8").Value = 1.742451
$ws.Range("O8").Value = 0.1705489461834183
$ws.Range("P8").Value = 0.1705489461834183
$ws.Range("Q8").Value = 2.274905304466667
$ws.Range("R8").Value = 20.4741477402
$ws.Range("S8").Value = 0.004447623599793239
$ws.Range("T8").Value = 0.004447623599793238

# Row 9
$ws.Range("G9").Value = 3.916733333333333
$ws.Range("H9").Value = 11.7502
$ws.Range("I9").Value = 0.0260782825067122
$ws.Range("J9").Value = 0.02607828250671219
$ws.Range("O9").Value = 0.6679715536912479
$ws.Range("P9").Value = 0.6679715536912479
$ws.Range("Q9").Value = 8.909888127311111
$ws.Range("R9").Value = 80.1889931458
$ws.Range("S9").Value = 0.01741955088360784
$ws.Range("T9").Value = 0.01741955088360784

# Row 10
$ws.Range("G10").Value = 3.916733333333333
$ws.Range("H10").Value = 11.7502
$ws.Range("I10").Value = 0.0260782825067122
$ws.Range("J10").Value = 0.02607828250671219
$ws.Range("M10").Value = 0.5499303333333333
$ws.Range("N10").Value = 1.649791
$ws.Range("O10").Value = 0.1614795001253337
$ws.Range("P10").Value = 0.1614795001253338
$ws.Range("Q10").Value = 2.153930467577778
$ws.Range("R10").Value = 19.3853742082
$ws.Range("S10").Value = 0.004211108023311121
$ws.Range("T10").Value = 0.004211108023311121

# Row 11
$ws.Range("G11").Value = 28.14764533333333
$ws.Range("H11").Value = 84.442936
$ws.Range("I11").Value = 0.1874118517730947
$ws.Range("J11").Value = 0.1874118517730947
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.580817
$ws.Range("N11").Value = 1.742451
$ws.Range("O11").Value = 0.1705489461834183
$ws.Range("P11").Value = 0.1705489461834183
$ws.Range("Q11").Value = 16.34863091957067
$ws.Range("R11").Value = 147.137678276136
$ws.Range("S11").Value = 0.03196289382218431
$ws.Range("T11").Value = 0.03196289382218431

# Row 12
$ws.Range("G12").Value = 28.14764533333333
$ws.Range("H12").Value = 84.442936
$ws.Range("I12").Value = 0.1874118517730947
$ws.Range("J12").Value = 0.1874118517730947
$ws.Range("O12").Value = 0.6679715536912479
$ws.Range("P12").Value = 0.6679715536912479
$ws.Range("Q12").Value = 64.03100482559378
$ws.Range("R12").Value = 576.279043430344
$ws.Range("S12").Value = 0.1251857858090279
$ws.Range("T12").Value = 0.1251857858090279

# Row 13
$ws.Range("G13").Value = 28.14764533333333
$ws.Range("H13").Value = 84.442936
$ws.Range("I13").Value = 0.1874118517730947
$ws.Range("J13").Value = 0.1874118517730947
$ws.Range("M13").Value = 0.5499303333333333
$ws.Range("N13").Value = 1.649791
$ws.Range("O13").Value = 0.1614795001253337
$ws.Range("P13").Value = 0.1614795001253338
$ws.Range("Q13").Value = 15.47924398070844
$ws.Range("R13").Value = 139.313195826376
$ws.Range("S13").Value = 0.03026317214188248
$ws.Range("T13").Value = 0.03026317214188248
